{"js": "// Add five new paragraphs of content at the end of the document. The\n// document currently ends with a single empty paragraph, so that\n// paragraph's text is replaced with the first new block, and four more\n// paragraphs are inserted after it for the remaining blocks.\n\nconst texts = [\n  \"Video provides a powerful way to help you prove your point. When you click Online Video, you can paste in the embed code for the video you want to add. You can also type a keyword to search online for the video that best fits your document.\",\n  \"To make your document look professionally produced, Word provides header, footer, cover page, and text box designs that complement each other. For example, you can add a matching cover page, header, and sidebar. Click Insert and then choose the elements you want from the different galleries.\",\n  \"Themes and styles also help keep your document coordinated. When you click Design and choose a new Theme, the pictures, charts, and SmartArt graphics change to match your new theme. When you apply styles, your headings change to match the new theme.\",\n  \"Save time in Word with new buttons that show up where you need them. To change the way a picture fits in your document, click it and a button for layout options appears next to it. When you work on a table, click where you want to add a row or a column, and then click the plus sign.\",\n  \"Reading is easier, too, in the new Reading view. You can collapse parts of the document and focus on the text you want. If you need to stop reading before you reach the end, Word remembers where you left off - even on another device.\"\n];\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// The trailing empty paragraph becomes the first new paragraph of text.\nlet current = body.paragraphs.items[body.paragraphs.items.length - 1];\ncurrent.insertText(texts[0], Word.InsertLocation.replace);\n\n// Each subsequent block is appended as a new paragraph after the previous one.\nfor (let i = 1; i < texts.length; i++) {\n  current = current.insertParagraph(texts[i], Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Add five new paragraphs of content at the end of the document, reusing\n# the existing final (empty) paragraph for the first new paragraph, and\n# inserting new paragraphs after it for the remaining four.\n\n$d = $word.ActiveDocument\n\n$texts = @(\n  \"Video provides a powerful way to help you prove your point. When you click Online Video, you can paste in the embed code for the video you want to add. You can also type a keyword to search online for the video that best fits your document.\",\n  \"To make your document look professionally produced, Word provides header, footer, cover page, and text box designs that complement each other. For example, you can add a matching cover page, header, and sidebar. Click Insert and then choose the elements you want from the different galleries.\",\n  \"Themes and styles also help keep your document coordinated. When you click Design and choose a new Theme, the pictures, charts, and SmartArt graphics change to match your new theme. When you apply styles, your headings change to match the new theme.\",\n  \"Save time in Word with new buttons that show up where you need them. To change the way a picture fits in your document, click it and a button for layout options appears next to it. When you work on a table, click where you want to add a row or a column, and then click the plus sign.\",\n  \"Reading is easier, too, in the new Reading view. You can collapse parts of the document and focus on the text you want. If you need to stop reading before you reach the end, Word remembers where you left off - even on another device.\"\n)\n\n# The document currently ends with a single empty paragraph; reuse it for\n# the first block of new text instead of leaving a stray blank paragraph.\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.Text = $texts[0]\n\nfor ($i = 1; $i -lt $texts.Count; $i++) {\n  $cur = $d.Paragraphs.Last\n  $cur.Range.InsertParagraphAfter()\n  $cur = $d.Paragraphs.Last\n  $cur.Range.Text = $texts[$i]\n}\n"}
